$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5592766404151917
$ws.Range("B1").Value = 3.960298776626587
$ws.Range("C1").Value = 6.227667331695557
$ws.Range("D1").Value = 1.469992876052856
$ws.Range("E1").Value = 0.844560980796814
